# Auto-generated edit script applying the Ultros_Profits commit diff.
# Updates currentAveragePrice / LevePrice / LeveProfit figures across
# several worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 550
$ws.Range("I20").Value = 550
$ws.Range("K20").Value = 550
$ws.Range("M20").Value = -320
$ws.Range("H35").Value = 550
$ws.Range("I35").Value = 550
$ws.Range("K35").Value = 550
$ws.Range("M35").Value = -171
$ws.Range("H92").Value = 348.45456
$ws.Range("I92").Value = 226.9375
$ws.Range("K92").Value = 226.9375
$ws.Range("M92").Value = 1021.0625
$ws.Range("H97").Value = 934.6
$ws.Range("J97").Value = 934.6
$ws.Range("L97").Value = 2803.8
$ws.Range("N97").Value = -3795.8
$ws.Range("H98").Value = 1250.909
$ws.Range("I98").Value = 875.4
$ws.Range("K98").Value = 875.4
$ws.Range("M98").Value = 622.6
$ws.Range("H99").Value = 449.33334
$ws.Range("I99").Value = 303
$ws.Range("J99").Value = 1181
$ws.Range("K99").Value = 909
$ws.Range("L99").Value = 3543
$ws.Range("M99").Value = 589
$ws.Range("N99").Value = -6539
$ws.Range("H112").Value = 3249.2856
$ws.Range("J112").Value = 3440.8333
$ws.Range("L112").Value = 10322.4999
$ws.Range("N112").Value = -12538.4999
$ws.Range("H122").Value = 1250.909
$ws.Range("I122").Value = 875.4
$ws.Range("K122").Value = 2626.2
$ws.Range("M122").Value = -176.1999999999998
$ws.Range("H127").Value = 4559.067
$ws.Range("I127").Value = 1860.5454
$ws.Range("J127").Value = 11980
$ws.Range("K127").Value = 5581.6362
$ws.Range("L127").Value = 35940
$ws.Range("M127").Value = -621.6361999999999
$ws.Range("N127").Value = -45860
$ws.Range("H133").Value = 57000
$ws.Range("J133").Value = 58888.89
$ws.Range("L133").Value = 58888.89
$ws.Range("N133").Value = -69008.89
$ws.Range("H137").Value = 3303.303
$ws.Range("I137").Value = 2315.5
$ws.Range("K137").Value = 6946.5
$ws.Range("M137").Value = -4396.5
$ws.Range("H141").Value = 9701.25
$ws.Range("I141").Value = 9046.111
$ws.Range("K141").Value = 27138.333
$ws.Range("M141").Value = -21958.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6273.7
$ws.Range("I74").Value = 5967.125
$ws.Range("K74").Value = 5967.125
$ws.Range("M74").Value = -5093.125
$ws.Range("H77").Value = 6273.7
$ws.Range("I77").Value = 5967.125
$ws.Range("K77").Value = 29835.625
$ws.Range("M77").Value = -25467.625
$ws.Range("H97").Value = 2716.7646
$ws.Range("I97").Value = 1555
$ws.Range("J97").Value = 4846.6665
$ws.Range("K97").Value = 1555
$ws.Range("L97").Value = 4846.6665
$ws.Range("M97").Value = -1059
$ws.Range("N97").Value = -5838.6665
$ws.Range("H132").Value = 2178.4614
$ws.Range("I132").Value = 1821.0555
$ws.Range("K132").Value = 5463.166499999999
$ws.Range("M132").Value = -2933.166499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 21
$ws.Range("I33").Value = 21
$ws.Range("K33").Value = 21
$ws.Range("M33").Value = 315
$ws.Range("H94").Value = 3832.923
$ws.Range("I94").Value = 916
$ws.Range("J94").Value = 8500
$ws.Range("K94").Value = 916
$ws.Range("L94").Value = 8500
$ws.Range("M94").Value = -465
$ws.Range("N94").Value = -9402
$ws.Range("H107").Value = 12502
$ws.Range("I107").Value = 11751.833
$ws.Range("K107").Value = 11751.833
$ws.Range("M107").Value = -9831.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1255.625
$ws.Range("I16").Value = 1091.5385
$ws.Range("K16").Value = 1091.5385
$ws.Range("M16").Value = -804.5385000000001
$ws.Range("H31").Value = 3518.05
$ws.Range("I31").Value = 2880.182
$ws.Range("J31").Value = 4297.6665
$ws.Range("K31").Value = 2880.182
$ws.Range("L31").Value = 4297.6665
$ws.Range("M31").Value = -2585.182
$ws.Range("N31").Value = -4887.6665
$ws.Range("H34").Value = 3518.05
$ws.Range("I34").Value = 2880.182
$ws.Range("J34").Value = 4297.6665
$ws.Range("K34").Value = 2880.182
$ws.Range("L34").Value = 4297.6665
$ws.Range("M34").Value = -2678.182
$ws.Range("N34").Value = -4701.6665
$ws.Range("H113").Value = 1255.625
$ws.Range("I113").Value = 1091.5385
$ws.Range("K113").Value = 1091.5385
$ws.Range("M113").Value = 1078.4615
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920
$ws.Range("H134").Value = 1418
$ws.Range("I134").Value = 1418
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4254
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1719
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 544.46875
$ws.Range("J12").Value = 596.6
$ws.Range("L12").Value = 1789.8
$ws.Range("N12").Value = -2135.8
$ws.Range("H23").Value = 2612.8572
$ws.Range("I23").Value = 5333
$ws.Range("J23").Value = 1871
$ws.Range("K23").Value = 15999
$ws.Range("L23").Value = 5613
$ws.Range("M23").Value = -15764
$ws.Range("N23").Value = -6083
$ws.Range("H107").Value = 1875.5476
$ws.Range("J107").Value = 1766.9678
$ws.Range("L107").Value = 5300.903399999999
$ws.Range("N107").Value = -9140.9034
$ws.Range("H112").Value = 4000
$ws.Range("I112").Value = 4000
$ws.Range("K112").Value = 12000
$ws.Range("M112").Value = -10892
$ws.Range("H122").Value = 6545.3076
$ws.Range("I122").Value = 272.25
$ws.Range("J122").Value = 9333.333
$ws.Range("K122").Value = 2450.25
$ws.Range("L122").Value = 83999.997
$ws.Range("M122").Value = -0.25
$ws.Range("N122").Value = -88899.997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 75780.75
$ws.Range("J124").Value = 75780.75
$ws.Range("L124").Value = 75780.75
$ws.Range("N124").Value = -85600.75
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 6520.216
$ws.Range("I132").Value = 6030.925
$ws.Range("K132").Value = 18092.775
$ws.Range("M132").Value = -15562.775

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H82").Value = 62502130
$ws.Range("J82").Value = 1853.1666
$ws.Range("L82").Value = 1853.1666
$ws.Range("N82").Value = -2575.1666
$ws.Range("H85").Value = 62502130
$ws.Range("J85").Value = 1853.1666
$ws.Range("L85").Value = 1853.1666
$ws.Range("N85").Value = -4349.1666
$ws.Range("H93").Value = 1323.3334
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13891626
$ws.Range("I81").Value = 2100
$ws.Range("J81").Value = 27781152
$ws.Range("K81").Value = 4200
$ws.Range("L81").Value = 55562304
$ws.Range("M81").Value = -3139
$ws.Range("N81").Value = -55564426
$ws.Range("H84").Value = 13891626
$ws.Range("I84").Value = 2100
$ws.Range("J84").Value = 27781152
$ws.Range("K84").Value = 21000
$ws.Range("L84").Value = 277811520
$ws.Range("M84").Value = -15696
$ws.Range("N84").Value = -277822128
$ws.Range("H122").Value = 2262
$ws.Range("I122").Value = 1924.6364
$ws.Range("K122").Value = 5773.9092
$ws.Range("M122").Value = -3323.9092
$ws.Range("H132").Value = 7141.2
$ws.Range("I132").Value = 3624.4614
$ws.Range("J132").Value = 30000
$ws.Range("K132").Value = 10873.3842
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -8343.3842
$ws.Range("N132").Value = -95060
$ws.Range("H136").Value = 70428.93
$ws.Range("I136").Value = 88186.82
$ws.Range("K136").Value = 264560.46
$ws.Range("M136").Value = -262010.46
